$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COMPOSIÇÃO DE VALORES")

$ws.Range("C3").Value = "Felipe Almeida"
$ws.Range("C4").Value = "Wellisson Chaves"
$ws.Range("C5").Value = "1020"

Write-Host "done"
